# Fixed names: singularize and pluralize
#
# The "Comments" sheet's "Constraints" column (C4, row for the "PostId"
# field) referenced the foreign table using the lower-case, un-pluralized
# form "posts". Fix it to match the actual sheet/table name "Posts".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comments")
$ws.Range("C4").Value = "required|fk:Posts,title"

# Leave the workbook with the Comments sheet active and C5 selected,
# matching where the author was working when the fix was saved.
$ws.Activate()
$ws.Range("C5").Select()
